# Auto-generated script to apply scheduled market-data refresh updates
# to the Belias_Profits workbook, per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3367.2703
$ws.Range("I64").Value = 3180.3462
$ws.Range("J64").Value = 3809.0908
$ws.Range("K64").Value = 3180.3462
$ws.Range("L64").Value = 3809.0908
$ws.Range("M64").Value = -2932.3462
$ws.Range("N64").Value = -4305.0908
$ws.Range("H67").Value = 3367.2703
$ws.Range("I67").Value = 3180.3462
$ws.Range("J67").Value = 3809.0908
$ws.Range("K67").Value = 3180.3462
$ws.Range("L67").Value = 3809.0908
$ws.Range("M67").Value = -2322.3462
$ws.Range("N67").Value = -5525.0908
$ws.Range("H113").Value = 3735.4285
$ws.Range("I113").Value = 2949.1428
$ws.Range("J113").Value = 4521.7144
$ws.Range("K113").Value = 2949.1428
$ws.Range("L113").Value = 4521.7144
$ws.Range("M113").Value = 304.8571999999999
$ws.Range("N113").Value = -11029.7144
$ws.Range("H137").Value = 1503076.9
$ws.Range("I137").Value = 1342.9
$ws.Range("J137").Value = 7939080
$ws.Range("K137").Value = 4028.7
$ws.Range("L137").Value = 23817240
$ws.Range("M137").Value = -1478.7
$ws.Range("N137").Value = -23822340

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 27477.785
$ws.Range("I74").Value = 31475.516
$ws.Range("J74").Value = 12819.444
$ws.Range("K74").Value = 31475.516
$ws.Range("L74").Value = 12819.444
$ws.Range("M74").Value = -30601.516
$ws.Range("N74").Value = -14567.444
$ws.Range("H77").Value = 27477.785
$ws.Range("I77").Value = 31475.516
$ws.Range("J77").Value = 12819.444
$ws.Range("K77").Value = 157377.58
$ws.Range("L77").Value = 64097.22
$ws.Range("M77").Value = -153009.58
$ws.Range("N77").Value = -72833.22
$ws.Range("H123").Value = 34333.332
$ws.Range("J123").Value = 34333.332
$ws.Range("L123").Value = 34333.332
$ws.Range("N123").Value = -44133.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1135.0555
$ws.Range("I80").Value = 1543.9
$ws.Range("J80").Value = 624
$ws.Range("K80").Value = 1543.9
$ws.Range("L80").Value = 624
$ws.Range("M80").Value = -545.9000000000001
$ws.Range("N80").Value = -2620
$ws.Range("H83").Value = 1135.0555
$ws.Range("I83").Value = 1543.9
$ws.Range("J83").Value = 624
$ws.Range("K83").Value = 7719.5
$ws.Range("L83").Value = 3120
$ws.Range("M83").Value = -2727.5
$ws.Range("N83").Value = -13104

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1052.079
$ws.Range("I31").Value = 936.5714
$ws.Range("J31").Value = 2399.6667
$ws.Range("K31").Value = 936.5714
$ws.Range("L31").Value = 2399.6667
$ws.Range("M31").Value = -641.5714
$ws.Range("N31").Value = -2989.6667
$ws.Range("H34").Value = 1052.079
$ws.Range("I34").Value = 936.5714
$ws.Range("J34").Value = 2399.6667
$ws.Range("K34").Value = 936.5714
$ws.Range("L34").Value = 2399.6667
$ws.Range("M34").Value = -734.5714
$ws.Range("N34").Value = -2803.6667
$ws.Range("H58").Value = 1870.0667
$ws.Range("I58").Value = 1920.9166
$ws.Range("J58").Value = 1666.6666
$ws.Range("K58").Value = 1920.9166
$ws.Range("L58").Value = 1666.6666
$ws.Range("M58").Value = -1717.9166
$ws.Range("N58").Value = -2072.6666
$ws.Range("H86").Value = 9816.286
$ws.Range("I86").Value = 4538.909
$ws.Range("J86").Value = 29166.666
$ws.Range("K86").Value = 4538.909
$ws.Range("L86").Value = 29166.666
$ws.Range("M86").Value = -3415.909
$ws.Range("N86").Value = -31412.666
$ws.Range("H89").Value = 9816.286
$ws.Range("I89").Value = 4538.909
$ws.Range("J89").Value = 29166.666
$ws.Range("K89").Value = 22694.545
$ws.Range("L89").Value = 145833.33
$ws.Range("M89").Value = -17078.545
$ws.Range("N89").Value = -157065.33
$ws.Range("H105").Value = 2005
$ws.Range("I105").Value = 2005
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2005
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -258
$ws.Range("N105").ClearContents()
$ws.Range("H134").Value = 2008.5416
$ws.Range("I134").Value = 2065
$ws.Range("J134").Value = 1794
$ws.Range("K134").Value = 6195
$ws.Range("L134").Value = 5382
$ws.Range("M134").Value = -3660
$ws.Range("N134").Value = -10452
$ws.Range("H136").Value = 1870.0667
$ws.Range("I136").Value = 1920.9166
$ws.Range("J136").Value = 1666.6666
$ws.Range("K136").Value = 5762.7498
$ws.Range("L136").Value = 4999.9998
$ws.Range("M136").Value = -3212.7498
$ws.Range("N136").Value = -10099.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 875.7
$ws.Range("I131").Value = 342.5
$ws.Range("J131").Value = 897.9167
$ws.Range("K131").Value = 1027.5
$ws.Range("L131").Value = 2693.7501
$ws.Range("M131").Value = 4012.5
$ws.Range("N131").Value = -12773.7501
$ws.Range("H140").Value = 2512.9412
$ws.Range("I140").Value = 837.1429000000001
$ws.Range("J140").Value = 10333.333
$ws.Range("K140").Value = 2511.4287
$ws.Range("L140").Value = 30999.999
$ws.Range("M140").Value = 2668.5713
$ws.Range("N140").Value = -41359.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4675.766
$ws.Range("I70").Value = 4445.9395
$ws.Range("J70").Value = 5217.5
$ws.Range("K70").Value = 4445.9395
$ws.Range("L70").Value = 5217.5
$ws.Range("M70").Value = -4175.9395
$ws.Range("N70").Value = -5757.5
$ws.Range("H73").Value = 4675.766
$ws.Range("I73").Value = 4445.9395
$ws.Range("J73").Value = 5217.5
$ws.Range("K73").Value = 4445.9395
$ws.Range("L73").Value = 5217.5
$ws.Range("M73").Value = -3509.9395
$ws.Range("N73").Value = -7089.5
$ws.Range("H107").Value = 730.53845
$ws.Range("I107").Value = 386.85715
$ws.Range("J107").Value = 1131.5
$ws.Range("K107").Value = 386.85715
$ws.Range("L107").Value = 1131.5
$ws.Range("M107").Value = 1533.14285
$ws.Range("N107").Value = -4971.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2036.5
$ws.Range("I7").Value = 1758.6
$ws.Range("K7").Value = 1758.6
$ws.Range("M7").Value = -1646.6
$ws.Range("H40").Value = 9360
$ws.Range("I40").Value = 9360
$ws.Range("K40").Value = 9360
$ws.Range("M40").Value = -9224
$ws.Range("H46").Value = 885.03845
$ws.Range("I46").Value = 1097.2858
$ws.Range("J46").Value = 806.8421
$ws.Range("K46").Value = 1097.2858
$ws.Range("L46").Value = 806.8421
$ws.Range("M46").Value = -909.2858000000001
$ws.Range("N46").Value = -1182.8421
$ws.Range("H55").Value = 680.8461
$ws.Range("I55").Value = 176.66667
$ws.Range("J55").Value = 832.1
$ws.Range("K55").Value = 176.66667
$ws.Range("L55").Value = 832.1
$ws.Range("M55").Value = -3.666670000000011
$ws.Range("N55").Value = -1178.1
$ws.Range("H126").Value = 2036.5
$ws.Range("I126").Value = 1758.6
$ws.Range("K126").Value = 5275.799999999999
$ws.Range("M126").Value = -2805.799999999999

Write-Output "Applied all market data updates."